$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "marker_1" column header (J1)
$ws.Range("J1").Value = "marker_1"

# Removing mislabelled wildtype strains (TDY451) from the perturbed
# samples: rows 5-7 carry genotype CNAG_06086 (a perturbation), so the
# wildtype "strain" label in column E no longer applies and is cleared.
$ws.Range("E5:E7").ClearContents()

$ws.Range("J6").Select() | Out-Null
